# Table Properties.xlsx edit
#
# Summary of the change (per commit message / diff):
#   - "Flags" column documented/stored as VARCHAR(3) instead of INTEGER.
#     This is reflected in the CREATE TABLE DDL string kept in cell H1 of
#     the "Properties Table" sheet.
#   - The per-row INSERT-OR-IGNORE helper formula in H2 now wraps the
#     Flags value (column C) in single quotes, matching the new VARCHAR
#     column type.
#   - The shared formula that used to populate H3:H11 with the same kind
#     of INSERT statement was removed entirely (cells cleared), leaving
#     only the first data row (H2) as a worked example.
#   - Minor selection / view-state touch-ups on a couple of sheets.

$wb = $excel.ActiveWorkbook

$wsField = $wb.Worksheets.Item("Field Explanation")
$wsTypes = $wb.Worksheets.Item("Types")
$wsProps = $wb.Worksheets.Item("Properties Table")

# --- "Properties Table" sheet -------------------------------------------------

# H1: the CREATE TABLE statement - Flags changes from INTEGER to VARCHAR(3)
$wsProps.Range("H1").Value = "CREATE TABLE [Properties] ([ID] INTEGER  NOT NULL PRIMARY KEY AUTOINCREMENT,[Name] VARCHAR(30)  UNIQUE NOT NULL,[Location] VARCHAR(150)  UNIQUE NOT NULL,[Flags] VARCHAR(3)  NOT NULL,[Staff] INTEGER DEFAULT '0' NOT NULL,[StaffCap] INTEGER DEFAULT '10' NOT NULL,[Cost] INTEGER  NOT NULL, [TypeID] INTEGER  NOT NULL)"

# H2: wrap the Flags value (C2) in quotes so it is emitted as a quoted
# string literal, matching the new VARCHAR(3) column type.
$wsProps.Range("H2").Formula = '="INSERT OR IGNORE INTO Properties (Name, Location, Flags, Staff, StaffCap, Cost, Type) VALUES (''" &A2&"'',''" &B2& "'',''"&C2&"'',"&D2&","&E2&","&F2&","&G2&");"'

# H3:H11 used to hold a shared copy of that formula (one INSERT statement
# per property row). Remove them completely - cell, formula and format.
$wsProps.Range("H3:H11").Clear()

# --- View / selection bookkeeping ---------------------------------------------

$wsField.Activate()
$wsField.Range("C9").Select()

# Leave "Types" sheet view state as-is (unchanged in the target).

# Re-activate "Properties Table" (keeps it the active/selected tab) and
# move the selection onto the freshly edited H2 cell.
$wsProps.Activate()
$wsProps.Range("H2").Select()
